$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.03850531578064
$ws.Range("B1").Value = 4.704395771026611
$ws.Range("C1").Value = 2.796947240829468
$ws.Range("D1").Value = 2.473300218582153
$ws.Range("E1").Value = 2.365476608276367
